$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 30.75612566666667
$ws.Cells.Item(2, 8).Value = 92.268377
$ws.Cells.Item(2, 9).Value = 0.9777985798685588
$ws.Cells.Item(2, 10).Value = 0.9777985798685588
$ws.Cells.Item(2, 13).Value = 4.877755666666666
$ws.Cells.Item(2, 14).Value = 14.633267
$ws.Cells.Item(2, 15).Value = 0.09961167132870688
$ws.Cells.Item(2, 16).Value = 0.09961167132870689
$ws.Cells.Item(2, 17).Value = 150.0208662552954
$ws.Cells.Item(2, 18).Value = 1350.187796297659
$ws.Cells.Item(2, 19).Value = 0.09740015076354322
$ws.Cells.Item(2, 20).Value = 0.09740015076354323
$ws.Cells.Item(3, 7).Value = 30.75612566666667
$ws.Cells.Item(3, 8).Value = 92.268377
$ws.Cells.Item(3, 9).Value = 0.9777985798685588
$ws.Cells.Item(3, 10).Value = 0.9777985798685588
$ws.Cells.Item(3, 15).Value = 0.1360673938501395
$ws.Cells.Item(3, 16).Value = 0.1360673938501395
$ws.Cells.Item(3, 17).Value = 204.9252665095643
$ws.Cells.Item(3, 18).Value = 1844.327398586079
$ws.Cells.Item(3, 19).Value = 0.1330465044730823
$ws.Cells.Item(3, 20).Value = 0.1330465044730823
$ws.Cells.Item(4, 7).Value = 30.75612566666667
$ws.Cells.Item(4, 8).Value = 92.268377
$ws.Cells.Item(4, 9).Value = 0.9777985798685588
$ws.Cells.Item(4, 10).Value = 0.9777985798685588
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.246459
$ws.Cells.Item(4, 14).Value = 0.739377
$ws.Cells.Item(4, 15).Value = 0.00503309197542868
$ws.Cells.Item(4, 16).Value = 0.00503309197542868
$ws.Cells.Item(4, 17).Value = 7.580123975680999
$ws.Cells.Item(4, 18).Value = 68.221115781129
$ws.Cells.Item(4, 19).Value = 0.004921350185922002
$ws.Cells.Item(4, 20).Value = 0.004921350185922002
$ws.Cells.Item(5, 7).Value = 30.75612566666667
$ws.Cells.Item(5, 8).Value = 92.268377
$ws.Cells.Item(5, 9).Value = 0.9777985798685588
$ws.Cells.Item(5, 10).Value = 0.9777985798685588
$ws.Cells.Item(5, 13).Value = 37.01331466666667
$ws.Cells.Item(5, 14).Value = 111.039944
$ws.Cells.Item(5, 15).Value = 0.7558718368280999
$ws.Cells.Item(5, 16).Value = 0.7558718368280999
$ws.Cells.Item(5, 17).Value = 1138.386157227876
$ws.Cells.Item(5, 18).Value = 10245.47541505089
$ws.Cells.Item(5, 19).Value = 0.739090408613155
$ws.Cells.Item(5, 20).Value = 0.739090408613155
$ws.Cells.Item(6, 7).Value = 30.75612566666667
$ws.Cells.Item(6, 8).Value = 92.268377
$ws.Cells.Item(6, 9).Value = 0.9777985798685588
$ws.Cells.Item(6, 10).Value = 0.9777985798685588
$ws.Cells.Item(6, 13).Value = 0.167274
$ws.Cells.Item(6, 14).Value = 0.501822
$ws.Cells.Item(6, 15).Value = 0.00341600601762507
$ws.Cells.Item(6, 16).Value = 0.00341600601762507
$ws.Cells.Item(6, 17).Value = 5.144700164766
$ws.Cells.Item(6, 18).Value = 46.302301482894
$ws.Cells.Item(6, 19).Value = 0.003340165832856244
$ws.Cells.Item(6, 20).Value = 0.003340165832856244
$ws.Cells.Item(7, 9).Value = 0.004830327290741966
$ws.Cells.Item(7, 10).Value = 0.004830327290741966
$ws.Cells.Item(7, 13).Value = 4.877755666666666
$ws.Cells.Item(7, 14).Value = 14.633267
$ws.Cells.Item(7, 15).Value = 0.09961167132870688
$ws.Cells.Item(7, 16).Value = 0.09961167132870689
$ws.Cells.Item(7, 17).Value = 0.7411034331335555
$ws.Cells.Item(7, 18).Value = 6.669930898202001
$ws.Cells.Item(7, 19).Value = 0.0004811569744954718
$ws.Cells.Item(7, 20).Value = 0.0004811569744954719
$ws.Cells.Item(8, 9).Value = 0.004830327290741966
$ws.Cells.Item(8, 10).Value = 0.004830327290741966
$ws.Cells.Item(8, 15).Value = 0.1360673938501395
$ws.Cells.Item(8, 16).Value = 0.1360673938501395
$ws.Cells.Item(8, 19).Value = 0.0006572500458944643
$ws.Cells.Item(8, 20).Value = 0.0006572500458944645
$ws.Cells.Item(9, 9).Value = 0.004830327290741966
$ws.Cells.Item(9, 10).Value = 0.004830327290741966
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.246459
$ws.Cells.Item(9, 14).Value = 0.739377
$ws.Cells.Item(9, 15).Value = 0.00503309197542868
$ws.Cells.Item(9, 16).Value = 0.00503309197542868
$ws.Cells.Item(9, 17).Value = 0.037445830318
$ws.Cells.Item(9, 18).Value = 0.337012472862
$ws.Cells.Item(9, 19).Value = 0.00002431148152572754
$ws.Cells.Item(9, 20).Value = 0.00002431148152572754
$ws.Cells.Item(10, 9).Value = 0.004830327290741966
$ws.Cells.Item(10, 10).Value = 0.004830327290741966
$ws.Cells.Item(10, 13).Value = 37.01331466666667
$ws.Cells.Item(10, 14).Value = 111.039944
$ws.Cells.Item(10, 15).Value = 0.7558718368280999
$ws.Cells.Item(10, 16).Value = 0.7558718368280999
$ws.Cells.Item(10, 17).Value = 5.623630301651556
$ws.Cells.Item(10, 18).Value = 50.612672714864
$ws.Cells.Item(10, 19).Value = 0.003651108361734029
$ws.Cells.Item(10, 20).Value = 0.003651108361734029
$ws.Cells.Item(11, 9).Value = 0.004830327290741966
$ws.Cells.Item(11, 10).Value = 0.004830327290741966
$ws.Cells.Item(11, 13).Value = 0.167274
$ws.Cells.Item(11, 14).Value = 0.501822
$ws.Cells.Item(11, 15).Value = 0.00341600601762507
$ws.Cells.Item(11, 16).Value = 0.00341600601762507
$ws.Cells.Item(11, 17).Value = 0.025414830948
$ws.Cells.Item(11, 18).Value = 0.228733478532
$ws.Cells.Item(11, 19).Value = 0.00001650042709227316
$ws.Cells.Item(11, 20).Value = 0.00001650042709227316
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.4673666666666667
$ws.Cells.Item(12, 8).Value = 1.4021
$ws.Cells.Item(12, 9).Value = 0.01485851852399773
$ws.Cells.Item(12, 10).Value = 0.01485851852399773
$ws.Cells.Item(12, 13).Value = 4.877755666666666
$ws.Cells.Item(12, 14).Value = 14.633267
$ws.Cells.Item(12, 15).Value = 0.09961167132870688
$ws.Cells.Item(12, 16).Value = 0.09961167132870689
$ws.Cells.Item(12, 17).Value = 2.279700406744444
$ws.Cells.Item(12, 18).Value = 20.5173036607
$ws.Cells.Item(12, 19).Value = 0.001480081863643965
$ws.Cells.Item(12, 20).Value = 0.001480081863643965
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.4673666666666667
$ws.Cells.Item(13, 8).Value = 1.4021
$ws.Cells.Item(13, 9).Value = 0.01485851852399773
$ws.Cells.Item(13, 10).Value = 0.01485851852399773
$ws.Cells.Item(13, 15).Value = 0.1360673938501395
$ws.Cells.Item(13, 16).Value = 0.1360673938501395
$ws.Cells.Item(13, 17).Value = 3.114021569633334
$ws.Cells.Item(13, 18).Value = 28.0261941267
$ws.Cells.Item(13, 19).Value = 0.002021759892034393
$ws.Cells.Item(13, 20).Value = 0.002021759892034394
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.4673666666666667
$ws.Cells.Item(14, 8).Value = 1.4021
$ws.Cells.Item(14, 9).Value = 0.01485851852399773
$ws.Cells.Item(14, 10).Value = 0.01485851852399773
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.246459
$ws.Cells.Item(14, 14).Value = 0.739377
$ws.Cells.Item(14, 15).Value = 0.00503309197542868
$ws.Cells.Item(14, 16).Value = 0.00503309197542868
$ws.Cells.Item(14, 17).Value = 0.1151867213
$ws.Cells.Item(14, 18).Value = 1.0366804917
$ws.Cells.Item(14, 19).Value = 0.00007478429034989138
$ws.Cells.Item(14, 20).Value = 0.00007478429034989138
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.4673666666666667
$ws.Cells.Item(15, 8).Value = 1.4021
$ws.Cells.Item(15, 9).Value = 0.01485851852399773
$ws.Cells.Item(15, 10).Value = 0.01485851852399773
$ws.Cells.Item(15, 13).Value = 37.01331466666667
$ws.Cells.Item(15, 14).Value = 111.039944
$ws.Cells.Item(15, 15).Value = 0.7558718368280999
$ws.Cells.Item(15, 16).Value = 0.7558718368280999
$ws.Cells.Item(15, 17).Value = 17.29878949804445
$ws.Cells.Item(15, 18).Value = 155.6891054824
$ws.Cells.Item(15, 19).Value = 0.01123113568927851
$ws.Cells.Item(15, 20).Value = 0.01123113568927851
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.4673666666666667
$ws.Cells.Item(16, 8).Value = 1.4021
$ws.Cells.Item(16, 9).Value = 0.01485851852399773
$ws.Cells.Item(16, 10).Value = 0.01485851852399773
$ws.Cells.Item(16, 13).Value = 0.167274
$ws.Cells.Item(16, 14).Value = 0.501822
$ws.Cells.Item(16, 15).Value = 0.00341600601762507
$ws.Cells.Item(16, 16).Value = 0.00341600601762507
$ws.Cells.Item(16, 17).Value = 0.07817829180000001
$ws.Cells.Item(16, 18).Value = 0.7036046262000001
$ws.Cells.Item(16, 19).Value = 0.00005075678869096983
$ws.Cells.Item(16, 20).Value = 0.00005075678869096983
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.07903166666666667
$ws.Cells.Item(17, 8).Value = 0.237095
$ws.Cells.Item(17, 9).Value = 0.002512574316701549
$ws.Cells.Item(17, 10).Value = 0.002512574316701549
$ws.Cells.Item(17, 13).Value = 4.877755666666666
$ws.Cells.Item(17, 14).Value = 14.633267
$ws.Cells.Item(17, 15).Value = 0.09961167132870688
$ws.Cells.Item(17, 16).Value = 0.09961167132870689
$ws.Cells.Item(17, 17).Value = 0.3854971599294444
$ws.Cells.Item(17, 18).Value = 3.469474439365
$ws.Cells.Item(17, 19).Value = 0.000250281727024225
$ws.Cells.Item(17, 20).Value = 0.000250281727024225
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 0.6666666666666666
$ws.Cells.Item(18, 7).Value = 0.07903166666666667
$ws.Cells.Item(18, 8).Value = 0.237095
$ws.Cells.Item(18, 9).Value = 0.002512574316701549
$ws.Cells.Item(18, 10).Value = 0.002512574316701549
$ws.Cells.Item(18, 15).Value = 0.1360673938501395
$ws.Cells.Item(18, 16).Value = 0.1360673938501395
$ws.Cells.Item(18, 17).Value = 0.5265808031183333
$ws.Cells.Item(18, 18).Value = 4.739227228065
$ws.Cells.Item(18, 19).Value = 0.0003418794391283748
$ws.Cells.Item(18, 20).Value = 0.0003418794391283749
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 0.6666666666666666
$ws.Cells.Item(19, 7).Value = 0.07903166666666667
$ws.Cells.Item(19, 8).Value = 0.237095
$ws.Cells.Item(19, 9).Value = 0.002512574316701549
$ws.Cells.Item(19, 10).Value = 0.002512574316701549
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.246459
$ws.Cells.Item(19, 14).Value = 0.739377
$ws.Cells.Item(19, 15).Value = 0.00503309197542868
$ws.Cells.Item(19, 16).Value = 0.00503309197542868
$ws.Cells.Item(19, 17).Value = 0.019478065535
$ws.Cells.Item(19, 18).Value = 0.175302589815
$ws.Cells.Item(19, 19).Value = 0.00001264601763105877
$ws.Cells.Item(19, 20).Value = 0.00001264601763105877
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 0.6666666666666666
$ws.Cells.Item(20, 7).Value = 0.07903166666666667
$ws.Cells.Item(20, 8).Value = 0.237095
$ws.Cells.Item(20, 9).Value = 0.002512574316701549
$ws.Cells.Item(20, 10).Value = 0.002512574316701549
$ws.Cells.Item(20, 13).Value = 37.01331466666667
$ws.Cells.Item(20, 14).Value = 111.039944
$ws.Cells.Item(20, 15).Value = 0.7558718368280999
$ws.Cells.Item(20, 16).Value = 0.7558718368280999
$ws.Cells.Item(20, 17).Value = 2.925223946964445
$ws.Cells.Item(20, 18).Value = 26.32701552268
$ws.Cells.Item(20, 19).Value = 0.001899184163932308
$ws.Cells.Item(20, 20).Value = 0.001899184163932308
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 0.6666666666666666
$ws.Cells.Item(21, 7).Value = 0.07903166666666667
$ws.Cells.Item(21, 8).Value = 0.237095
$ws.Cells.Item(21, 9).Value = 0.002512574316701549
$ws.Cells.Item(21, 10).Value = 0.002512574316701549
$ws.Cells.Item(21, 13).Value = 0.167274
$ws.Cells.Item(21, 14).Value = 0.501822
$ws.Cells.Item(21, 15).Value = 0.00341600601762507
$ws.Cells.Item(21, 16).Value = 0.00341600601762507
$ws.Cells.Item(21, 17).Value = 0.01321994301
$ws.Cells.Item(21, 18).Value = 0.11897948709
$ws.Cells.Item(21, 19).Value = 0.00000858296898558269
$ws.Cells.Item(21, 20).Value = 0.00000858296898558269
